$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row to the table "Condicion_Pacientes" (expands table ref/autoFilter and sheet dimension)
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $tbl.ListRows.Add()

# Copy formatting (incl. date number format) from the row above into the new row
$ws.Range("A42:F42").Copy()
$ws.Range("A43").PasteSpecial(-4122)

# Fill in the new data
$ws.Range("A43").Value = 43962
$ws.Range("B43").Value = 349
$ws.Range("C43").Value = 128
$ws.Range("D43").Value = 225
$ws.Range("E43").Value = 11
$ws.Range("F43").Value = 15

$ws.Range("D46").Select()
